$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete worker detail rows (WALFREDY JULIO MUENTES, 3 period rows)
$ws.Range("17:19").EntireRow.Delete()

# Refresh the summary figures for the updated account-statement data
$ws.Range("E11").Value = 128000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
